$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the extra barcode/invoice data, keeping only the location column values
$ws.Range("C2").ClearContents()
$ws.Range("B3").ClearContents()
$ws.Range("C3").ClearContents()

# Update selection to B2
$ws.Range("B2").Select()
